$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Preserve the header style (bold, centered, thin border) before rewriting the row
$ws.Range("A1").Copy()
$ws.Range("A1:BA1").PasteSpecial(-4122)

# Write the full header row (A1:BA1) with the new column layout
$ws.Range("A1").Value = "salutation [Link]"
$ws.Range("B1").Value = "customer_name [Data]"
$ws.Range("C1").Value = "customer_type [Select]"
$ws.Range("D1").Value = "customer_group [Link]"
$ws.Range("E1").Value = "territory [Link]"
$ws.Range("F1").Value = "gender [Link]"
$ws.Range("G1").Value = "lead_name [Link]"
$ws.Range("H1").Value = "opportunity_name [Link]"
$ws.Range("I1").Value = "prospect_name [Link]"
$ws.Range("J1").Value = "account_manager [Link]"
$ws.Range("K1").Value = "default_currency [Link]"
$ws.Range("L1").Value = "default_bank_account [Link]"
$ws.Range("M1").Value = "default_price_list [Link]"
$ws.Range("N1").Value = "is_internal_customer [Check]"
$ws.Range("O1").Value = "represents_company [Link]"
$ws.Range("P1").Value = "market_segment [Link]"
$ws.Range("Q1").Value = "industry [Link]"
$ws.Range("R1").Value = "website [Data]"
$ws.Range("S1").Value = "language [Link]"
$ws.Range("T1").Value = "customer_details [Text]"
$ws.Range("U1").Value = "customer_primary_address [Link]"
$ws.Range("V1").Value = "customer_primary_contact [Link]"
$ws.Range("W1").Value = "tax_id [Data]"
$ws.Range("X1").Value = "tax_category [Link]"
$ws.Range("Y1").Value = "tax_withholding_category [Link]"
$ws.Range("Z1").Value = "payment_terms [Link]"
$ws.Range("AA1").Value = "credit_limits.1.credit_limit [Currency]"
$ws.Range("AB1").Value = "credit_limits.1.bypass_credit_limit_check [Check]"
$ws.Range("AC1").Value = "credit_limits.2.credit_limit [Currency]"
$ws.Range("AD1").Value = "credit_limits.2.bypass_credit_limit_check [Check]"
$ws.Range("AE1").Value = "credit_limits.3.credit_limit [Currency]"
$ws.Range("AF1").Value = "credit_limits.3.bypass_credit_limit_check [Check]"
$ws.Range("AG1").Value = "credit_limits.4.credit_limit [Currency]"
$ws.Range("AH1").Value = "credit_limits.4.bypass_credit_limit_check [Check]"
$ws.Range("AI1").Value = "credit_limits.5.credit_limit [Currency]"
$ws.Range("AJ1").Value = "credit_limits.5.bypass_credit_limit_check [Check]"
$ws.Range("AK1").Value = "loyalty_program [Link]"
$ws.Range("AL1").Value = "sales_team.1.allocated_percentage [Float]"
$ws.Range("AM1").Value = "sales_team.1.incentives [Currency]"
$ws.Range("AN1").Value = "sales_team.2.allocated_percentage [Float]"
$ws.Range("AO1").Value = "sales_team.2.incentives [Currency]"
$ws.Range("AP1").Value = "sales_team.3.allocated_percentage [Float]"
$ws.Range("AQ1").Value = "sales_team.3.incentives [Currency]"
$ws.Range("AR1").Value = "sales_team.4.allocated_percentage [Float]"
$ws.Range("AS1").Value = "sales_team.4.incentives [Currency]"
$ws.Range("AT1").Value = "sales_team.5.allocated_percentage [Float]"
$ws.Range("AU1").Value = "sales_team.5.incentives [Currency]"
$ws.Range("AV1").Value = "default_sales_partner [Link]"
$ws.Range("AW1").Value = "default_commission_rate [Float]"
$ws.Range("AX1").Value = "so_required [Check]"
$ws.Range("AY1").Value = "dn_required [Check]"
$ws.Range("AZ1").Value = "is_frozen [Check]"
$ws.Range("BA1").Value = "disabled [Check]"

$ws.Application.CutCopyMode = $false
